$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 149
$ws.Cells.Item(149, 2).Value = 48654
$ws.Cells.Item(149, 5).Value = 38.26
$ws.Cells.Item(149, 6).Value = -1
$ws.Cells.Item(149, 7).Value = -32.02

# Row 150
$ws.Cells.Item(150, 2).Value = 63902
$ws.Cells.Item(150, 5).Value = 34.04
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 7).Value = 64.04000000000001

# Row 161
$ws.Cells.Item(161, 2).Value = 53925
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 66.44

# Row 162
$ws.Cells.Item(162, 2).Value = 64350
$ws.Cells.Item(162, 5).Value = 70.63
$ws.Cells.Item(162, 6).Value = 60
$ws.Cells.Item(162, 7).Value = 3986.4

# Row 163
$ws.Cells.Item(163, 2).Value = 57756
$ws.Cells.Item(163, 5).Value = 79.37
$ws.Cells.Item(163, 6).Value = -100
$ws.Cells.Item(163, 7).Value = -6644

# Row 183
$ws.Cells.Item(183, 2).Value = 57552
$ws.Cells.Item(183, 5).Value = 136.86
$ws.Cells.Item(183, 6).Value = -5
$ws.Cells.Item(183, 7).Value = -603.45

# Row 184
$ws.Cells.Item(184, 2).Value = 64329
$ws.Cells.Item(184, 5).Value = 128.32
$ws.Cells.Item(184, 6).Value = 6
$ws.Cells.Item(184, 7).Value = 724.14

# Row 317
$ws.Cells.Item(317, 2).Value = 61610
$ws.Cells.Item(317, 4).Value = 102.71
$ws.Cells.Item(317, 5).Value = 122.71
$ws.Cells.Item(317, 6).Value = -58
$ws.Cells.Item(317, 7).Value = -5957.18

# Row 318
$ws.Cells.Item(318, 2).Value = 57077
$ws.Cells.Item(318, 4).Value = 93.08
$ws.Cells.Item(318, 5).Value = 111.2
$ws.Cells.Item(318, 6).Value = 1
$ws.Cells.Item(318, 7).Value = 93.08

# Row 346
$ws.Cells.Item(346, 2).Value = 63520
$ws.Cells.Item(346, 5).Value = 153.4
$ws.Cells.Item(346, 6).Value = 94
$ws.Cells.Item(346, 7).Value = 13562.32

# Row 347
$ws.Cells.Item(347, 2).Value = 55373
$ws.Cells.Item(347, 5).Value = 163.62
$ws.Cells.Item(347, 6).Value = -94
$ws.Cells.Item(347, 7).Value = -13562.32

# Row 350
$ws.Cells.Item(350, 2).Value = 57802
$ws.Cells.Item(350, 5).Value = 162.71
$ws.Cells.Item(350, 6).Value = -79
$ws.Cells.Item(350, 7).Value = -11334.92

# Row 351
$ws.Cells.Item(351, 2).Value = 63571
$ws.Cells.Item(351, 6).Value = 19
$ws.Cells.Item(351, 7).Value = 2726.12

# Row 352
$ws.Cells.Item(352, 2).Value = 63531
$ws.Cells.Item(352, 5).Value = 152.53
$ws.Cells.Item(352, 6).Value = 80
$ws.Cells.Item(352, 7).Value = 11478.4

# Row 355
$ws.Cells.Item(355, 2).Value = 55356
$ws.Cells.Item(355, 5).Value = 54.04
$ws.Cells.Item(355, 6).Value = -158
$ws.Cells.Item(355, 7).Value = -7527.12

# Row 356
$ws.Cells.Item(356, 2).Value = 63510
$ws.Cells.Item(356, 5).Value = 50.66
$ws.Cells.Item(356, 6).Value = 159
$ws.Cells.Item(356, 7).Value = 7574.76

# Row 372
$ws.Cells.Item(372, 2).Value = 63652
$ws.Cells.Item(372, 5).Value = 55.42
$ws.Cells.Item(372, 6).Value = 223
$ws.Cells.Item(372, 7).Value = 11624.99

# Row 373
$ws.Cells.Item(373, 2).Value = 57885
$ws.Cells.Item(373, 5).Value = 62.28
$ws.Cells.Item(373, 6).Value = 4
$ws.Cells.Item(373, 7).Value = 208.52

# Row 379
$ws.Cells.Item(379, 2).Value = 63564
$ws.Cells.Item(379, 6).Value = 55
$ws.Cells.Item(379, 7).Value = 7095.55

# Row 380
$ws.Cells.Item(380, 2).Value = 65514
$ws.Cells.Item(380, 6).Value = 1
$ws.Cells.Item(380, 7).Value = 129.01

# Row 382
$ws.Cells.Item(382, 2).Value = 63560
$ws.Cells.Item(382, 5).Value = 134.87
$ws.Cells.Item(382, 6).Value = 52
$ws.Cells.Item(382, 7).Value = 6596.72

# Row 383
$ws.Cells.Item(383, 2).Value = 60325
$ws.Cells.Item(383, 5).Value = 151.57
$ws.Cells.Item(383, 6).Value = -102
$ws.Cells.Item(383, 7).Value = -12939.72

# Row 389
$ws.Cells.Item(389, 2).Value = 57817
$ws.Cells.Item(389, 6).Value = 3
$ws.Cells.Item(389, 7).Value = 239.43

# Row 390
$ws.Cells.Item(390, 2).Value = 62865
$ws.Cells.Item(390, 6).Value = 73
$ws.Cells.Item(390, 7).Value = 5826.13

# Row 419
$ws.Cells.Item(419, 2).Value = 63007
$ws.Cells.Item(419, 6).Value = 898
$ws.Cells.Item(419, 7).Value = 153854.34

# Row 420
$ws.Cells.Item(420, 2).Value = 57856
$ws.Cells.Item(420, 6).Value = 2
$ws.Cells.Item(420, 7).Value = 342.66

# Row 431
$ws.Cells.Item(431, 2).Value = 53082
$ws.Cells.Item(431, 3).Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Cells.Item(431, 6).Value = 1
$ws.Cells.Item(431, 7).Value = 59.47

# Row 432
$ws.Cells.Item(432, 2).Value = 63102
$ws.Cells.Item(432, 3).Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Cells.Item(432, 6).Value = 4
$ws.Cells.Item(432, 7).Value = 237.88

# Row 457
$ws.Cells.Item(457, 2).Value = 63681
$ws.Cells.Item(457, 5).Value = 23.84
$ws.Cells.Item(457, 6).Value = 43
$ws.Cells.Item(457, 7).Value = 964.0599999999999

# Row 458
$ws.Cells.Item(458, 2).Value = 31930
$ws.Cells.Item(458, 5).Value = 26.8
$ws.Cells.Item(458, 6).Value = -62
$ws.Cells.Item(458, 7).Value = -1390.04

# Row 536
$ws.Cells.Item(536, 2).Value = 47097
$ws.Cells.Item(536, 4).Value = 112.28
$ws.Cells.Item(536, 5).Value = 134.16
$ws.Cells.Item(536, 6).Value = 15
$ws.Cells.Item(536, 7).Value = 1684.2

# Row 537
$ws.Cells.Item(537, 2).Value = 58047
$ws.Cells.Item(537, 4).Value = 105.54
$ws.Cells.Item(537, 5).Value = 126.1
$ws.Cells.Item(537, 6).Value = 51
$ws.Cells.Item(537, 7).Value = 5382.54

# Row 579
$ws.Cells.Item(579, 2).Value = 65069
$ws.Cells.Item(579, 5).Value = 14.3
$ws.Cells.Item(579, 6).Value = 121
$ws.Cells.Item(579, 7).Value = 1627.45

# Row 580
$ws.Cells.Item(580, 2).Value = 53757
$ws.Cells.Item(580, 5).Value = 16.08
$ws.Cells.Item(580, 6).Value = -159
$ws.Cells.Item(580, 7).Value = -2138.55

# Row 590
$ws.Cells.Item(590, 2).Value = 64922
$ws.Cells.Item(590, 5).Value = 20.98
$ws.Cells.Item(590, 6).Value = 176
$ws.Cells.Item(590, 7).Value = 3472.48

# Row 591
$ws.Cells.Item(591, 2).Value = 45706
$ws.Cells.Item(591, 5).Value = 23.58
$ws.Cells.Item(591, 6).Value = -202
$ws.Cells.Item(591, 7).Value = -3985.46

# Row 593
$ws.Cells.Item(593, 2).Value = 45718
$ws.Cells.Item(593, 5).Value = 19.38
$ws.Cells.Item(593, 6).Value = -294
$ws.Cells.Item(593, 7).Value = -4768.68

# Row 594
$ws.Cells.Item(594, 2).Value = 64927
$ws.Cells.Item(594, 5).Value = 17.26
$ws.Cells.Item(594, 6).Value = 286
$ws.Cells.Item(594, 7).Value = 4638.92

# Row 709
$ws.Cells.Item(709, 2).Value = 64833
$ws.Cells.Item(709, 5).Value = 34.9
$ws.Cells.Item(709, 6).Value = 99
$ws.Cells.Item(709, 7).Value = 3250.17

# Row 710
$ws.Cells.Item(710, 2).Value = 60025
$ws.Cells.Item(710, 5).Value = 37.22
$ws.Cells.Item(710, 6).Value = -98
$ws.Cells.Item(710, 7).Value = -3217.34

# Row 720
$ws.Cells.Item(720, 2).Value = 64830
$ws.Cells.Item(720, 5).Value = 34.9
$ws.Cells.Item(720, 6).Value = 117
$ws.Cells.Item(720, 7).Value = 3841.11

# Row 721
$ws.Cells.Item(721, 2).Value = 60022
$ws.Cells.Item(721, 5).Value = 37.22
$ws.Cells.Item(721, 6).Value = -113
$ws.Cells.Item(721, 7).Value = -3709.79

# Row 859
$ws.Cells.Item(859, 2).Value = 63150
$ws.Cells.Item(859, 4).Value = 75.68000000000001
$ws.Cells.Item(859, 5).Value = 80.45
$ws.Cells.Item(859, 6).Value = 193
$ws.Cells.Item(859, 7).Value = 14606.24

# Row 860
$ws.Cells.Item(860, 2).Value = 61428
$ws.Cells.Item(860, 4).Value = 69.16
$ws.Cells.Item(860, 5).Value = 73.52
$ws.Cells.Item(860, 6).Value = 1
$ws.Cells.Item(860, 7).Value = 69.16
